# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 0
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
